$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("G:G").Insert()
